$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'28.124.54"
$ws.Range("E2").Value = "  -1.52%  "
# Row 3
$ws.Range("D3").Value = "'1.895.74"
$ws.Range("E3").Value = "  -0.87%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'314.53"
$ws.Range("E5").Value = "  -0.05%  "
# Row 6
$ws.Range("E6").Value = "  +0.00%  "
# Row 7
$ws.Range("D7").Value = "'0.5024"
$ws.Range("E7").Value = "  -0.65%  "
# Row 8
$ws.Range("D8").Value = "'0.3901"
$ws.Range("E8").Value = "  -1.63%  "
# Row 9
$ws.Range("D9").Value = "'0.09233"
$ws.Range("E9").Value = "  -5.63%  "
# Row 10
$ws.Range("D10").Value = "'1.130"
$ws.Range("E10").Value = "  -2.53%  "
# Row 11
$ws.Range("D11").Value = "'41.89"
$ws.Range("E11").Value = "  +0.19%  "
# Row 12
$ws.Range("D12").Value = "'6.397"
$ws.Range("E12").Value = "  -2.32%  "
# Row 13
$ws.Range("E13").Value = "  -1.66%  "
# Row 14
$ws.Range("D14").Value = "'1.905.84"
$ws.Range("E14").Value = "  -0.12%  "
# Row 15
$ws.Range("E15").Value = "  -3.79%  "
# Row 16
$ws.Range("E16").Value = "  +0.05%  "
# Row 17
$ws.Range("D17").Value = "'0.00001112"
$ws.Range("E17").Value = "  -2.60%  "
# Row 18
$ws.Range("D18").Value = "'92.48"
$ws.Range("E18").Value = "  -1.50%  "
# Row 19
$ws.Range("D19").Value = "'0.06653"
$ws.Range("E19").Value = "  -0.05%  "
# Row 20
$ws.Range("E20").Value = "  -1.08%  "
# Row 21
$ws.Range("E21").Value = "  -0.01%  "
# Row 22
$ws.Range("D22").Value = "'6.213"
$ws.Range("E22").Value = "  -1.35%  "
# Row 23
$ws.Range("D23").Value = "'28.188.82"
$ws.Range("E23").Value = "  -1.50%  "
# Row 24
$ws.Range("E24").Value = "  +0.01%  "
# Row 25
$ws.Range("D25").Value = "'2.320"
$ws.Range("E25").Value = "  +1.83%  "
# Row 26
$ws.Range("D26").Value = "'2.126.00"
$ws.Range("E26").Value = "  -0.06%  "
# Row 27
$ws.Range("D27").Value = "'2.556"
$ws.Range("E27").Value = "  -7.24%  "
# Row 28
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").Value = "'158.44"
$ws.Range("E28").Value = "  -0.62%  "
# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'20.87"
$ws.Range("E29").Value = "  -2.12%  "
# Row 30
$ws.Range("D30").Value = "'127.02"
$ws.Range("E30").Value = "  -1.55%  "
# Row 31
$ws.Range("E31").Value = "  -2.28%  "
# Row 32
$ws.Range("D32").Value = "'0.1060"
$ws.Range("E32").Value = "  -1.29%  "
# Row 33
$ws.Range("D33").Value = "'5.609"
$ws.Range("E33").Value = "  -2.35%  "
# Row 34
$ws.Range("D34").Value = "'3.623"
$ws.Range("E34").Value = "  -0.51%  "
# Row 35
$ws.Range("D35").Value = "'9.566"
$ws.Range("E35").Value = "  -3.37%  "
# Row 36
$ws.Range("D36").Value = "'0.06596"
$ws.Range("E36").Value = "  -3.18%  "
# Row 37
$ws.Range("D37").Value = "'1.341"
$ws.Range("E37").Value = "  +12.53%  "
# Row 38
$ws.Range("D38").Value = "'0.02400"
$ws.Range("E38").Value = "  -1.95%  "
# Row 39
$ws.Range("D39").Value = "'0.2204"
$ws.Range("E39").Value = "  -1.41%  "
# Row 40
$ws.Range("E40").Value = "  -3.94%  "
# Row 41
$ws.Range("D41").Value = "'0.6474"
$ws.Range("E41").Value = "  +0.47%  "
# Row 42
$ws.Range("D42").Value = "'4.977"
$ws.Range("E42").Value = "  -2.66%  "
# Row 43
$ws.Range("D43").Value = "'11.41"
$ws.Range("E43").Value = "  -2.86%  "
# Row 44
$ws.Range("E44").Value = "  -0.03%  "
# Row 45
$ws.Range("D45").Value = "'0.6106"
$ws.Range("E45").Value = "  -0.04%  "
# Row 46
$ws.Range("D46").Value = "'13.38"
$ws.Range("E46").Value = "  -2.27%  "
# Row 47
$ws.Range("E47").Value = "  +1.42%  "
# Row 48
$ws.Range("D48").Value = "'3.692"
$ws.Range("E48").Value = "  +0.56%  "
# Row 49
$ws.Range("D49").Value = "'2.003"
$ws.Range("E49").Value = "  -2.10%  "
# Row 50
$ws.Range("D50").Value = "'122.20"
# Row 51
$ws.Range("D51").Value = "'1.205"
$ws.Range("E51").Value = "  -0.79%  "
